# Fix Arbeidsgiverbelop bug in "Multiple Employers" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Multiple Employers")

# --- Root-cause input corrections -------------------------------------
# Grad (percentage sick) for Emp 2 was 10%, should be 80%.
$ws.Range("E8").Value = 0.8
# Refusjonsgrad for Emp 2 was 10%, should be 90%.
$ws.Range("E10").Value = 0.9

# --- Situation narrative text (row 32) ---------------------------------
$ws.Range("C32").Formula = '=IF(SUM(D30:G31)<=MAKSBELOP,"Everyone paid",IF(SUM(D30:G30)<=MAKSBELOP,"Arbeidsgivere fully paid; Person partially paid by Person request", "Arbeidsgivere partial payment ratio by Arbeidsgiver request"))'

# --- Arbeidsgiverbelop scaling fix (row 33) -----------------------------
# Bug: used ScalingFactor (6G based) instead of the MAKSBELOP ratio when
# employers alone exceed MAKSBELOP.
# NOTE: formulas below use single quotes so that PowerShell does not try
# to interpolate "$D"/"$G" etc. as variables.
$ws.Range("D33").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,D30,IF(SUM($D$30:$G$30)<=MAKSBELOP,D30,ROUND(D30*MAKSBELOP/SUM($D$30:$G$30),0)))'
$ws.Range("E33").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,E30,IF(SUM($D$30:$G$30)<=MAKSBELOP,E30,ROUND(E30*MAKSBELOP/SUM($D$30:$G$30),0)))'
$ws.Range("F33").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,F30,IF(SUM($D$30:$G$30)<=MAKSBELOP,F30,ROUND(F30*MAKSBELOP/SUM($D$30:$G$30),0)))'
$ws.Range("G33").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,G30,IF(SUM($D$30:$G$30)<=MAKSBELOP,G30,ROUND(G30*MAKSBELOP/SUM($D$30:$G$30),0)))'

# Row 34 (PersonRemainder helper row) is no longer interesting on its own
# now that row 36 supersedes it as the rounding-error check; hide it.
$ws.Rows.Item(34).Hidden = $true

# --- Personbelop fix (row 35) -------------------------------------------
# Only spread PersonRemainder across persons when employers alone are
# within MAKSBELOP; otherwise persons get nothing (or full amount if
# everyone fits).
$ws.Range("D35").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,D31,IF(SUM($D$30:$G$30)<=MAKSBELOP,ROUND(PersonRemainder*D31/SUM($D$31:$G$31),0),0))'
$ws.Range("E35").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,E31,IF(SUM($D$30:$G$30)<=MAKSBELOP,ROUND(PersonRemainder*E31/SUM($D$31:$G$31),0),0))'
$ws.Range("F35").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,F31,IF(SUM($D$30:$G$30)<=MAKSBELOP,ROUND(PersonRemainder*F31/SUM($D$31:$G$31),0),0))'
$ws.Range("G35").Formula = '=IF(SUM($D$30:$G$31)<=MAKSBELOP,G31,IF(SUM($D$30:$G$30)<=MAKSBELOP,ROUND(PersonRemainder*G31/SUM($D$31:$G$31),0),0))'

# --- New row 36: overall daily rounding error check ---------------------
$ws.Range("A34:G34").Copy()
$ws.Range("A36:G36").PasteSpecial(-4122)
$ws.Range("H35").Copy()
$ws.Range("H36").PasteSpecial(-4122)

$ws.Range("B36").Value = "Daily rounding error"
$ws.Range("C36").Formula = "=MAKSBELOP-SUM(D33:G35)"
$ws.Range("H36").Value = "Integer"

# --- Cosmetic view state (best effort) ----------------------------------
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.Zoom = 140
$ws.Range("D3").Select()
